$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.266.24'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').Value = '1.592.92'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '213.14'
$ws.Range('E5').Value = '  +0.72%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.499'
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.52%  '
$ws.Range('E9').Value = '  -0.61%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '18.98'
$ws.Range('E10').Value = '  -2.29%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0850'
$ws.Range('E11').Value = '  +0.35%  '
$ws.Range('D12').Value = '1.817.12'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('D13').Value = '1.594.17'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.00'
$ws.Range('E14').Value = '  -1.06%  '
$ws.Range('E15').Value = '  -2.46%  '
$ws.Range('D17').Value = '26.247.75'
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('D18').Value = '0.0₃0722'
$ws.Range('E18').Value = '  -1.40%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '215.79'
$ws.Range('E19').Value = '  +1.51%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.38'
$ws.Range('E20').Value = '  -1.65%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.29'
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('E23').Value = '  +0.28%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.10'
$ws.Range('E24').Value = '  -3.18%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '144.73'
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  -1.38%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.113'
$ws.Range('E28').Value = '  +0.60%  '
$ws.Range('E29').Value = '  -0.70%  '
$ws.Range('E30').Value = '  -2.33%  '
$ws.Range('E31').Value = '  +0.69%  '
$ws.Range('E32').Value = '  -0.64%  '
$ws.Range('D33').Value = '1.407.85'
$ws.Range('E33').Value = '  +5.18%  '
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('E35').Value = '  -0.74%  '
$ws.Range('E36').Value = '  -1.89%  '
$ws.Range('E37').Value = '  -4.44%  '
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.823'
$ws.Range('E39').Value = '  +0.31%  '
$ws.Range('E40').Value = '  -0.36%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.964'
$ws.Range('E42').Value = '  -7.80%  '
$ws.Range('E43').Value = '  +0.78%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.760'
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('D45').Value = '1.729.22'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '60.81'
$ws.Range('E46').Value = '  -1.63%  '
$ws.Range('E47').Value = '  -0.75%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.49'
$ws.Range('E48').Value = '  -1.56%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0502'
$ws.Range('E49').Value = '  -0.55%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0950'
$ws.Range('E50').Value = '  -3.39%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.14%  '
